# src: routes: - updated the csv import
#
# The template's "department" column instructs the CSV importer how to
# separate multiple department names. It used to use "/" as the
# separator; switch it to "|" instead. This string is shared by the
# three sample rows (F2:F4) in the template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F4").Value = "name of department 1|name of department 2"

# Restore the sheet's last-used selection/active cell, as left by the
# author after editing (previously F4).
$ws.Range("E8").Select()
